$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure sheet ---
# Insert two new columns (D, E) for "ownTeam" and "oppTeam" before the existing
# "batsman" column (which shifts from D to F, and the rest shift accordingly).
$ws.Columns("D:E").Insert()

# Insert a new row above the current row 2 to hold the "Dubai (DSC) Oct 24 2020"
# match, which becomes the new first data row.
$ws.Rows("2:2").Insert()

# Rows 5 and 6 (for the remaining new matches) don't exist yet, so they will
# simply be created when values are written into them below.

# --- Header row ---
$ws.Range("A1").Value = "venue"
$ws.Range("B1").Value = "date"
$ws.Range("C1").Value = "result"
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"
$ws.Range("F1").Value = "batsman"
$ws.Range("G1").Value = "totalRuns"
$ws.Range("H1").Value = "totalBalls"
$ws.Range("I1").Value = "total4s"
$ws.Range("J1").Value = "total6s"
$ws.Range("K1").Value = "sr"

# Force the numeric-looking columns to be stored as text, matching the
# source data (numbers stored as text).
$ws.Range("G2:K6").NumberFormat = "@"

# --- Row 2: Dubai (DSC), October 24 2020 vs Sunrisers Hyderabad ---
$ws.Range("A2").Value = " Dubai (DSC)"
$ws.Range("B2").Value = " October 24 2020"
$ws.Range("C2").Value = "Kings XI won by 12 runs"
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Sunrisers Hyderabad"
$ws.Range("F2").Value = "Chris Jordan "
$ws.Range("G2").Value = "7"
$ws.Range("H2").Value = "12"
$ws.Range("I2").Value = "0"
$ws.Range("J2").Value = "0"
$ws.Range("K2").Value = "58.33"

# --- Row 3: Abu Dhabi, November 01 2020 vs Chennai Super Kings ---
$ws.Range("A3").Value = " Abu Dhabi"
$ws.Range("B3").Value = " November 01 2020"
$ws.Range("C3").Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Range("D3").Value = "Kings XI Punjab"
$ws.Range("E3").Value = "Chennai Super Kings"
$ws.Range("F3").Value = "Chris Jordan "
$ws.Range("G3").Value = "4"
$ws.Range("H3").Value = "5"
$ws.Range("I3").Value = "0"
$ws.Range("J3").Value = "0"
$ws.Range("K3").Value = "80.00"

# --- Row 4: Abu Dhabi, October 10 2020 vs Kolkata Knight Riders ---
$ws.Range("A4").Value = " Abu Dhabi"
$ws.Range("B4").Value = " October 10 2020"
$ws.Range("C4").Value = "KKR won by 2 runs"
$ws.Range("D4").Value = "Kings XI Punjab"
$ws.Range("E4").Value = "Kolkata Knight Riders"
$ws.Range("F4").Value = "Chris Jordan "
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "-"

# --- Row 5: Dubai (DSC), September 20 2020 vs Delhi Capitals ---
$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " September 20 2020"
$ws.Range("C5").Value = "Match tied (Capitals won the one-over eliminator)"
$ws.Range("D5").Value = "Kings XI Punjab"
$ws.Range("E5").Value = "Delhi Capitals"
$ws.Range("F5").Value = "Chris Jordan "
$ws.Range("G5").Value = "5"
$ws.Range("H5").Value = "6"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "83.33"

# --- Row 6: Dubai (DSC), October 18 2020 vs Mumbai Indians ---
$ws.Range("A6").Value = " Dubai (DSC)"
$ws.Range("B6").Value = " October 18 2020"
$ws.Range("C6").Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Range("D6").Value = "Kings XI Punjab"
$ws.Range("E6").Value = "Mumbai Indians"
$ws.Range("F6").Value = "Chris Jordan "
$ws.Range("G6").Value = "13"
$ws.Range("H6").Value = "8"
$ws.Range("I6").Value = "2"
$ws.Range("J6").Value = "0"
$ws.Range("K6").Value = "162.50"
